$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVERAGES")

# Update the Percent value in C2 to the value that used to be in C3
$ws.Range("C2").Value = 24.22

# Remove row 3 entirely (it was a duplicate divider row), shrinking the used range to A1:C2
$ws.Rows(3).Delete()
